$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date (column C) for rows 2 through 16
# from 2023-10-13 (serial 45212) to 2023-10-22 (serial 45221).
$ws.Range("C2:C16").Value = 45221
